$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from the last existing header cell (H1) onto the
# two new header cells, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data cells in row 2 (plain, unstyled, numeric - matching H2 etc.)
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
